$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCellValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextCellValue $ws.Range('D2') '51.765.92'
Set-TextCellValue $ws.Range('E2') '  -0.42%  '
Set-TextCellValue $ws.Range('D3') '2.786.46'
Set-TextCellValue $ws.Range('E3') '  -0.14%  '
Set-TextCellValue $ws.Range('E4') '  +0.06%  '
Set-TextCellValue $ws.Range('D5') '353.38'
Set-TextCellValue $ws.Range('E5') '  -1.47%  '
Set-TextCellValue $ws.Range('D6') '108.88'
Set-TextCellValue $ws.Range('E6') '  -0.56%  '
Set-TextCellValue $ws.Range('D7') '0.548'
Set-TextCellValue $ws.Range('E7') '  -3.31%  '
Set-TextCellValue $ws.Range('E8') '  +0.36%  '
Set-TextCellValue $ws.Range('D9') '0.596'
Set-TextCellValue $ws.Range('E9') '  +0.08%  '
Set-TextCellValue $ws.Range('D10') '39.93'
Set-TextCellValue $ws.Range('E10') '  -0.54%  '
Set-TextCellValue $ws.Range('E11') '  +2.92%  '
Set-TextCellValue $ws.Range('E12') '  +3.64%  '
Set-TextCellValue $ws.Range('D13') '0.0838'
Set-TextCellValue $ws.Range('E13') '  -2.14%  '
Set-TextCellValue $ws.Range('E14') '  +0.29%  '
Set-TextCellValue $ws.Range('D15') '3.231.03'
Set-TextCellValue $ws.Range('E15') '  +0.13%  '
Set-TextCellValue $ws.Range('D16') '2.767.34'
Set-TextCellValue $ws.Range('E16') '  -1.26%  '
Set-TextCellValue $ws.Range('D17') '0.928'
Set-TextCellValue $ws.Range('E17') '  -1.44%  '
Set-TextCellValue $ws.Range('D18') '51.760.14'
Set-TextCellValue $ws.Range('E18') '  -0.27%  '
Set-TextCellValue $ws.Range('E19') '  +4.30%  '
Set-TextCellValue $ws.Range('E20') '  -0.30%  '
Set-TextCellValue $ws.Range('D21') '13.16'
Set-TextCellValue $ws.Range('E21') '  +0.54%  '
Set-TextCellValue $ws.Range('E22') '  -1.68%  '
Set-TextCellValue $ws.Range('D23') '70.00'
Set-TextCellValue $ws.Range('E23') '  -0.41%  '
Set-TextCellValue $ws.Range('D24') '266.61'
Set-TextCellValue $ws.Range('E24') '  -2.94%  '
Set-TextCellValue $ws.Range('E25') '  -0.33%  '
Set-TextCellValue $ws.Range('D26') '26.15'
Set-TextCellValue $ws.Range('E26') '  -2.16%  '
Set-TextCellValue $ws.Range('D27') '0.999'
Set-TextCellValue $ws.Range('E27') '  -0.14%  '
Set-TextCellValue $ws.Range('E28') '  +11.17%  '
Set-TextCellValue $ws.Range('E29') '  +0.35%  '
Set-TextCellValue $ws.Range('D30') '36.87'
Set-TextCellValue $ws.Range('E30') '  +7.43%  '
Set-TextCellValue $ws.Range('E31') '  +8.65%  '
Set-TextCellValue $ws.Range('D32') '51.87'
Set-TextCellValue $ws.Range('E32') '  +0.58%  '
Set-TextCellValue $ws.Range('D33') '0.0454'
Set-TextCellValue $ws.Range('E33') '  -2.30%  '
Set-TextCellValue $ws.Range('E34') '  +6.32%  '
Set-TextCellValue $ws.Range('D35') '2.07'
Set-TextCellValue $ws.Range('E35') '  -8.50%  '
Set-TextCellValue $ws.Range('E36') '  -1.87%  '
Set-TextCellValue $ws.Range('E37') '  +0.13%  '
Set-TextCellValue $ws.Range('D38') '18.51'
Set-TextCellValue $ws.Range('E38') '  +2.53%  '
Set-TextCellValue $ws.Range('D39') '3.15'
Set-TextCellValue $ws.Range('E39') '  -2.91%  '
Set-TextCellValue $ws.Range('E40') '  -1.70%  '
Set-TextCellValue $ws.Range('D41') '2.54'
Set-TextCellValue $ws.Range('E41') '  -1.30%  '
Set-TextCellValue $ws.Range('E42') '  -0.83%  '
Set-TextCellValue $ws.Range('D43') '120.89'
Set-TextCellValue $ws.Range('E43') '  -1.20%  '
Set-TextCellValue $ws.Range('D44') '21.99'
Set-TextCellValue $ws.Range('E44') '  -0.14%  '
Set-TextCellValue $ws.Range('E45') '  -2.52%  '
Set-TextCellValue $ws.Range('D46') '2.131.49'
Set-TextCellValue $ws.Range('E46') '  +2.65%  '
Set-TextCellValue $ws.Range('D47') '3.29'
Set-TextCellValue $ws.Range('E47') '  +0.91%  '
Set-TextCellValue $ws.Range('D48') '2.34'
Set-TextCellValue $ws.Range('E48') '  +7.44%  '
Set-TextCellValue $ws.Range('E49') '  -5.25%  '
Set-TextCellValue $ws.Range('D50') '0.907'
Set-TextCellValue $ws.Range('E50') '  -3.58%  '
Set-TextCellValue $ws.Range('D51') '1.34'
Set-TextCellValue $ws.Range('E51') '  +9.33%  '
